$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the row for Unit 5 ("0/1") entirely ---
$ws.Rows.Item(5).Delete()

# --- Unit 12 and Unit 14 changed from "0/?" to "0/2" ---
$ws.Cells.Item(7, 2).Value = "0/2"
$ws.Cells.Item(8, 2).Value = "0/2"

# --- Remove the rows for Unit 17, Unit 18 and Unit 26 ---
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()

# --- The row that used to be Unit 40 ("1/1") becomes Unit 31 ("1/2") ---
$ws.Cells.Item(12, 1).Value = 31
$ws.Cells.Item(12, 2).Value = "1/2"

# --- Unit 28 changed from "0.4/1" to "0.4/2" ---
$ws.Cells.Item(9, 2).Value = "0.4/2"

# --- The row that used to be Unit 42 ("0/1") becomes Unit 40 ("2/2") ---
$ws.Cells.Item(13, 1).Value = 40
$ws.Cells.Item(13, 2).Value = "2/2"

# --- Leave a styled (text-formatted) but empty cell at B17 ---
$ws.Cells.Item(17, 2).NumberFormat = "@"

# --- Update the active selection ---
$ws.Range("D12").Select()
